# Adding "Cargo" (Job) and "Departamento" (Department) columns to the
# student/user import template, after the existing "Años Experiencia" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells - these also append the two new shared strings
# ("Cargo", "Departamento") and extend the sheet dimension to A1:I1.
$ws.Range("H1").Value = "Cargo"
$ws.Range("I1").Value = "Departamento"

# Give the new columns a sensible width, matching the "best fit" sizing
# that Excel applies to the other header columns.
$ws.Columns.Item(1).ColumnWidth = 7.6640625
$ws.Columns.Item(2).ColumnWidth = 6.6640625
$ws.Columns.Item(3).ColumnWidth = 8.33203125
$ws.Columns.Item(4).ColumnWidth = 4.21875
$ws.Columns.Item(5).ColumnWidth = 28.88671875
$ws.Columns.Item(6).ColumnWidth = 25.33203125
$ws.Columns.Item(7).ColumnWidth = 14.77734375

# Move the active selection to where the author left off editing (I2,
# right below the newly typed "Departamento" header).
$ws.Range("I2").Select()
